$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "47.296.78"
$ws.Range("E2").Value = "  +2.48%  "
$ws.Range("D3").Value = "2.504.93"
$ws.Range("E3").Value = "  +2.25%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "324.28"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.68%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "109.20"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +4.40%  "
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("E9").Value = "  +0.18%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.15"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +8.87%  "
$ws.Range("E11").Value = "  +1.03%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.123"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.67%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "18.41"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.76%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.19"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.96%  "
$ws.Range("D15").Value = "2.897.63"
$ws.Range("E15").Value = "  +2.30%  "
$ws.Range("D16").Value = "2.504.71"
$ws.Range("E16").Value = "  +1.73%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.857"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.82%  "
$ws.Range("D18").Value = "47.228.46"
$ws.Range("E18").Value = "  +2.66%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.86"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.72%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.65"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.80%  "
$ws.Range("E21").Value = "  +1.07%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.71"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +13.74%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "70.50"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.59%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "247.65"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.43%  "
$ws.Range("E25").Value = "  +3.07%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "26.02"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.44%  "
$ws.Range("E27").Value = "  -0.08%  "
$ws.Range("E28").Value = "  +0.75%  "
$ws.Range("E29").Value = "  +3.85%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "35.61"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +4.40%  "
$ws.Range("E31").Value = "  +8.30%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "20.04"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.04%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.45"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.19%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0793"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.99%  "
$ws.Range("E36").Value = "  +0.32%  "
$ws.Range("E37").Value = "  +5.05%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.71"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +4.26%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.00"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.59%  "
$ws.Range("E40").Value = "  +1.21%  "
$ws.Range("E41").Value = "  +0.25%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "121.38"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -5.08%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "21.23"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.17%  "
$ws.Range("E44").Value = "  +2.45%  "
$ws.Range("D45").Value = "1.996.75"
$ws.Range("E45").Value = "  +1.48%  "
$ws.Range("E46").Value = "  +4.17%  "
$ws.Range("E47").Value = "  -1.12%  "
$ws.Range("B48").Value = "FraxShare"
$ws.Range("C48").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.09"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.01%  "
$ws.Range("B49").Value = "Stacks"
$ws.Range("C49").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.78"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -4.68%  "
$ws.Range("E50").Value = "  +3.88%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "56.90"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +4.58%  "
